$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new ACH receiver rows (18: Republic Services, 19: New Pig) ---
# Values are entered in an order that matches how the new shared strings
# were appended to the workbook (names first, then the routing/account
# numbers that needed to stay text so leading zeros / exact digit strings
# are preserved).

$ws.Range("A18").Value = "Republic"
$ws.Range("B18").Value = "Republic Services"

$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "121000248"

$ws.Range("A19").Value = "New Pig"
$ws.Range("B19").Value = "New Pig"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "362198604"

$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "5053406"

$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "613027"

$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "031301422"

# D18 (account number) stays a real number, but keeps the "text" style (s=2)
# used throughout column D - set the value first, then apply the format so
# the stored type remains numeric.
$ws.Range("D18").Value = 4140909680
$ws.Range("D18").NumberFormat = "@"

$ws.Range("E18").Value = "vendor"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("G18").Value = "checking"

$ws.Range("E19").Value = "vendor"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("G19").Value = "checking"

# --- Minor view/formatting touch-ups ---
$ws.Rows(6).RowHeight = 18.55

$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 205
$ws.Range("C20").Select() | Out-Null
